$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the three changed labels in the pin-mapping table
$ws.Range("C2").Value = "Manuca 5V"
$ws.Range("H2").Value = "Manuca 3.3V"
$ws.Range("C3").Value = "CN301.32 (PA4)"

# Match the saved selection from the commit (cursor left on C2)
$ws.Range("C2").Select()
